$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "batsman" column (D), shifting the
# existing D:I ("batsman".."sr") over to F:K.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells.
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# New data for the existing row 2.
$ws.Range("D2").Value = "Chennai Super Kings"
$ws.Range("E2").Value = "Sunrisers Hyderabad"

# New row 3: another Dwayne Bravo batting line.
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 10 2020"
$ws.Range("C3").Value = "RCB won by 37 runs"
$ws.Range("D3").Value = "Chennai Super Kings"
$ws.Range("E3").Value = "Royal Challengers Bangalore"

# Copy the batsman name from the cell it already lives in (after the column
# insert above, "Dwayne Bravo " now sits in F2) so the exact characters
# (incl. the trailing non-breaking space) are preserved.
$ws.Range("F2").Copy($ws.Range("F3"))

# Numeric-looking stats need to stay text (matches the rest of the sheet,
# which stores every value as text), so force Text format before assigning.
$statsRow3 = $ws.Range("G3:K3")
$statsRow3.NumberFormat = "@"
$ws.Range("G3").Value = "7"
$ws.Range("H3").Value = "5"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "140.00"
